$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Se agregan los scripts 0147-0155 a la clase Tests_AdmInstituciones
# The old row 64 (blank B/C placeholder) is pushed down by the 9 new
# data rows, so insert 9 fresh rows right above it and fill them in
# with the same repeating pattern used by the preceding test-case rows.
$newCodes = @("DEC_0147","DEC_0148","DEC_0149","DEC_0150","DEC_0151","DEC_0152","DEC_0153","DEC_0154","DEC_0155")

for ($i = 0; $i -lt $newCodes.Length; $i++) {
    $r = 64 + $i
    $ws.Rows.Item($r).Insert()
    $ws.Cells.Item($r, 1).Value = $newCodes[$i]
    $ws.Cells.Item($r, 2).Value = "13712759-8"
    $ws.Cells.Item($r, 3).Value = "Verity1.1"
    $ws.Cells.Item($r, 4).Value = "SIN_DATO"
    $ws.Cells.Item($r, 5).Value = "SIN_DATO"
    $ws.Cells.Item($r, 6).Value = "SIN_DATO"
    $ws.Cells.Item($r, 7).Value = "SIN_DATO"
    $ws.Cells.Item($r, 8).Value = "SIN_DATO"
    $ws.Cells.Item($r, 9).Value = "SIN_DATO"
    $ws.Cells.Item($r, 10).Value = "SIN_DATO"
}

# Two more blank spacer rows (matching the blank B/C placeholder row that
# used to sit right above the trailing summary rows) got inserted as well.
$ws.Rows.Item(74).Insert()
$ws.Rows.Item(74).Insert()

# Reflect where the user ended up after adding the rows.
$ws.Range("C63").Select()
